# Apply updated Betfair Back/Lay odds values for 2025-10-09 (rows 2-7, columns F:AO)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.34
$ws.Range("G2").Value = 2.4
$ws.Range("H2").Value = 3.45
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 3.4
$ws.Range("K2").Value = 3.5
$ws.Range("L2").Value = 1.51
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 3.15
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 1.72
$ws.Range("Q2").Value = 2.34
$ws.Range("R2").Value = 1.25
$ws.Range("S2").Value = 4.6
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 1.98
$ws.Range("V2").Value = 1.38
$ws.Range("W2").Value = 1.72
$ws.Range("X2").Value = 11.5
$ws.Range("Y2").Value = 11.5
$ws.Range("Z2").Value = 24
$ws.Range("AA2").Value = 75
$ws.Range("AC2").Value = 7.4
$ws.Range("AD2").Value = 15
$ws.Range("AE2").Value = 48
$ws.Range("AF2").Value = 14
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 95
$ws.Range("AJ2").Value = 110
$ws.Range("AK2").Value = 28
$ws.Range("AL2").Value = 130
$ws.Range("AM2").Value = 140
$ws.Range("AN2").Value = 38
$ws.Range("AO2").Value = 60

# Row 3
$ws.Range("F3").Value = 1.88
$ws.Range("G3").Value = 1.89
$ws.Range("H3").Value = 5.5
$ws.Range("I3").Value = 5.7
$ws.Range("J3").Value = 3.45
$ws.Range("L3").Value = 1.58
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 2.92
$ws.Range("O3").Value = 1.51
$ws.Range("P3").Value = 1.6
$ws.Range("Q3").Value = 2.58
$ws.Range("R3").Value = 1.21
$ws.Range("S3").Value = 5.3
$ws.Range("T3").Value = 2.22
$ws.Range("U3").Value = 1.73
$ws.Range("V3").Value = 1.21
$ws.Range("W3").Value = 2.12
$ws.Range("X3").Value = 9.6
$ws.Range("Y3").Value = 14.5
$ws.Range("AB3").Value = 6.8
$ws.Range("AC3").Value = 7.8
$ws.Range("AF3").Value = 9.4
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 29
$ws.Range("AI3").Value = 130
$ws.Range("AJ3").Value = 22
$ws.Range("AK3").Value = 26
$ws.Range("AL3").Value = 60
$ws.Range("AM3").Value = 330
$ws.Range("AN3").Value = 22
$ws.Range("AO3").Value = 170

# Row 4
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2.04
$ws.Range("H4").Value = 4.4
$ws.Range("I4").Value = 4.7
$ws.Range("K4").Value = 3.55
$ws.Range("L4").Value = 1.49
$ws.Range("O4").Value = 1.41
$ws.Range("P4").Value = 1.78
$ws.Range("Q4").Value = 2.22
$ws.Range("R4").Value = 1.28
$ws.Range("S4").Value = 4.2
$ws.Range("U4").Value = 1.94
$ws.Range("V4").Value = 1.28
$ws.Range("X4").Value = 12
$ws.Range("Z4").Value = 32
$ws.Range("AC4").Value = 7.8
$ws.Range("AD4").Value = 18
$ws.Range("AE4").Value = 65
$ws.Range("AF4").Value = 11.5
$ws.Range("AJ4").Value = 25
$ws.Range("AK4").Value = 24
$ws.Range("AL4").Value = 44
$ws.Range("AN4").Value = 18
$ws.Range("AO4").Value = 80

# Row 5
$ws.Range("F5").Value = 3.5
$ws.Range("G5").Value = 3.6
$ws.Range("H5").Value = 2.3
$ws.Range("I5").Value = 2.36
$ws.Range("L5").Value = 1.49
$ws.Range("M5").Value = 1.09
$ws.Range("N5").Value = 3.15
$ws.Range("O5").Value = 1.41
$ws.Range("P5").Value = 1.74
$ws.Range("Q5").Value = 2.22
$ws.Range("R5").Value = 1.27
$ws.Range("S5").Value = 4.2
$ws.Range("T5").Value = 1.92
$ws.Range("U5").Value = 1.9
$ws.Range("V5").Value = 1.73
$ws.Range("W5").Value = 1.38
$ws.Range("X5").Value = 11.5
$ws.Range("Y5").Value = 8.800000000000001
$ws.Range("Z5").Value = 13.5
$ws.Range("AB5").Value = 12
$ws.Range("AF5").Value = 23
$ws.Range("AH5").Value = 20
$ws.Range("AJ5").Value = 70
$ws.Range("AK5").Value = 46
$ws.Range("AL5").Value = 65
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 25

# Row 6
$ws.Range("F6").Value = 3.95
$ws.Range("G6").Value = 4.2
$ws.Range("H6").Value = 2.22
$ws.Range("I6").Value = 2.26
$ws.Range("J6").Value = 3.2
$ws.Range("K6").Value = 3.3
$ws.Range("L6").Value = 1.56
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 2.74
$ws.Range("O6").Value = 1.54
$ws.Range("P6").Value = 1.57
$ws.Range("Q6").Value = 2.64
$ws.Range("R6").Value = 1.2
$ws.Range("S6").Value = 5.4
$ws.Range("T6").Value = 2.16
$ws.Range("U6").Value = 1.77
$ws.Range("V6").Value = 1.79
$ws.Range("W6").Value = 1.32
$ws.Range("X6").Value = 9
$ws.Range("Y6").Value = 7.4
$ws.Range("Z6").Value = 12
$ws.Range("AA6").Value = 960
$ws.Range("AB6").Value = 11
$ws.Range("AC6").Value = 8
$ws.Range("AD6").Value = 12
$ws.Range("AE6").Value = 30
$ws.Range("AF6").Value = 27
$ws.Range("AG6").Value = 23
$ws.Range("AH6").Value = 25
$ws.Range("AI6").Value = 960
$ws.Range("AJ6").Value = 100
$ws.Range("AK6").Value = 70
$ws.Range("AL6").Value = 1000
$ws.Range("AM6").Value = 180
$ws.Range("AO6").Value = 32

# Row 7
$ws.Range("F7").Value = 3.05
$ws.Range("G7").Value = 3.25
$ws.Range("H7").Value = 2.7
$ws.Range("I7").Value = 2.76
$ws.Range("J7").Value = 3.05
$ws.Range("K7").Value = 3.2
$ws.Range("P7").Value = 1.67
$ws.Range("U7").Value = 1.9
$ws.Range("V7").Value = 1.57
$ws.Range("W7").Value = 1.44
$ws.Range("X7").Value = 11.5
$ws.Range("AC7").Value = 42
$ws.Range("AH7").Value = 980
